$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 42
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 27
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 11
$ws.Range("G2").Value = 37
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 22
$ws.Range("B3").Value = 2
$ws.Range("D3").Value = 7
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 4
$ws.Range("B4").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 1
$ws.Range("J5").Value = 2
$ws.Range("D8").Value = 7
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 5
$ws.Range("C12").Value = 1
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 4
$ws.Range("D13").Value = 1
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 1
$ws.Range("B17").Value = 2
$ws.Range("D17").Value = 9
$ws.Range("G17").Value = 1
$ws.Range("B19").Value = 2
$ws.Range("D19").Value = 9
$ws.Range("G19").Value = 1
$ws.Range("B23").Value = 8
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 8
$ws.Range("J23").Value = 3
$ws.Range("C24").Value = 1
$ws.Range("H24").Value = 1
$ws.Range("B25").Value = 1
$ws.Range("D25").Value = 1
$ws.Range("H25").Value = 1
$ws.Range("J25").Value = 3
$ws.Range("B26").Value = 7
$ws.Range("E26").Value = 2
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 5
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 1
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 3
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 5
$ws.Range("I28").Value = 3
$ws.Range("J28").Value = 8
$ws.Range("D29").Value = 2
$ws.Range("H29").Value = 2
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 6
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 1
$ws.Range("J30").Value = 1
$ws.Range("G31").Value = 1
$ws.Range("C32").Value = 4
$ws.Range("H32").Value = 2
$ws.Range("J32").Value = 1
$ws.Range("H33").Value = 1
$ws.Range("I33").Value = 2
$ws.Range("B34").Value = 13
$ws.Range("C34").Value = 3
$ws.Range("D34").Value = 3
$ws.Range("E34").Value = 3
$ws.Range("F34").Value = 8
$ws.Range("G34").Value = 7
$ws.Range("H34").Value = 4
$ws.Range("I34").Value = 1
$ws.Range("J34").Value = 1
$ws.Range("B35").Value = 6
$ws.Range("C35").Value = 3
$ws.Range("D35").Value = 3
$ws.Range("G35").Value = 1
$ws.Range("H36").Value = 0
$ws.Range("B37").Value = 4
$ws.Range("E37").Value = 3
$ws.Range("F37").Value = 8
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 4
$ws.Range("I37").Value = 1
$ws.Range("J37").Value = 1
$ws.Range("B38").Value = 3
$ws.Range("B40").Value = 15
$ws.Range("C40").Value = 1
$ws.Range("D40").Value = 3
$ws.Range("G40").Value = 17
$ws.Range("H40").Value = 4
$ws.Range("J40").Value = 1
$ws.Range("B41").Value = 6
$ws.Range("G41").Value = 15
$ws.Range("H41").Value = 1
$ws.Range("B42").Value = 5
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 3
$ws.Range("H42").Value = 1
$ws.Range("J42").Value = 1
$ws.Range("B43").Value = 1
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 2
$ws.Range("B44").Value = 3
